$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the blank "placeholder" cells (AR15, AR16) from the existing
# blank AR14 cell before AR14 itself gets cleared below.
$ws.Range("AR14").Copy($ws.Range("AR15"))
$ws.Range("AR14").Copy($ws.Range("AR16"))

# --- Row 14 ---
$ws.Range("A14").Value = 111798760
$ws.Range("Q14").Value = 753108.8301749222
$ws.Range("R14").Value = 7091007.708399305
$ws.Range("S14").Value = 100
$ws.Range("AR14").ClearContents()

# --- Row 15 ---
$ws.Range("A15").Value = 111798757
$ws.Range("B15").Value = 81076
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5046
$ws.Range("F15").Value = "Grön jordtunga"
$ws.Range("G15").Value = "Microglossum viride"
$ws.Range("H15").Value = "(Pers.:Fr.) Gillet"
$ws.Range("AF15").Value = "mikroskoperad"
$ws.Range("AI15").Value = "Granskog"

# --- Row 16 ---
$ws.Range("A16").Value = 111798755
$ws.Range("B16").Value = 90709
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5448
$ws.Range("F16").Value = "Svartvit taggsvamp"
$ws.Range("G16").Value = "Phellodon connatus"
$ws.Range("H16").Value = "(Schultz) nom.prov"
$ws.Range("Q16").Value = 753030.7189070459
$ws.Range("R16").Value = 7090920.781295684
$ws.Range("S16").Value = 25
$ws.Range("AF16").ClearContents()
